$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ATATATATAT"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("A3").Value = "ATATATAAAT"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "7"
$ws.Range("G3").Value = ""
$ws.Range("A4").Value = "AAATATATAT"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "1"
$ws.Range("G4").Value = ""
$ws.Range("A5").Value = "TAATATATAT"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "1"
$ws.Range("A6").Value = "ATATAAATAT"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "5"
$ws.Range("G6").Value = ""
$ws.Range("A7").Value = "ATAAATATAT"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "3"
$ws.Range("G7").Value = ""
$ws.Range("A8").Value = "TAAAATATAT"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "3"
$ws.Range("G8").Value = "1"
$ws.Range("A9").Value = "TATAATATAT"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = "1, 3"
$ws.Range("A10").Value = "ATAAATAAAT"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = "3, 7"
$ws.Range("G10").Value = ""
$ws.Range("A11").Value = "AATAATATAT"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = "1"
$ws.Range("G11").Value = "3"
$ws.Range("A12").Value = "ATATAATAAT"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = "5"
$ws.Range("G12").Value = "7"
$ws.Range("A13").Value = "AAATAAATAT"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = "1, 5"
$ws.Range("G13").Value = ""
$ws.Range("A14").Value = "TAATAAATAT"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = "5"
$ws.Range("G14").Value = "1"
$ws.Range("A15").Value = "ATATAAAAAT"
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "5, 7"
$ws.Range("G15").Value = ""
$ws.Range("A16").Value = "AAATATAAAT"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = "1, 7"
$ws.Range("G16").Value = ""
$ws.Range("A17").Value = "AAAAATATAT"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = "1, 3"
$ws.Range("G17").Value = ""
$ws.Range("A18").Value = "TAATATAAAT"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = "7"
$ws.Range("G18").Value = "1"
$ws.Range("A19").Value = "TAATAAAAAT"
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = "5, 7"
$ws.Range("G19").Value = "1"
$ws.Range("A20").Value = "AAAAATAAAT"
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = "1, 3, 7"
$ws.Range("G20").Value = ""
$ws.Range("A21").Value = "AAATAAAAAT"
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = "1, 5, 7"
$ws.Range("G21").Value = ""
$ws.Range("A22").Value = "AATAATAAAT"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = "1, 7"
$ws.Range("G22").Value = "3"
$ws.Range("A23").Value = "TAAAATAAAT"
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = "3, 7"
$ws.Range("G23").Value = "1"
$ws.Range("A24").Value = "AAATAATAAT"
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = "1, 5"
$ws.Range("G24").Value = "7"
$ws.Range("A25").Value = "TATAATAAAT"
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = "7"
$ws.Range("G25").Value = "1, 3"